$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New consolidated values for A2:A7
$ws.Range("A2").Value = "('森', ['Basic Land — Forest', '({T}: Add {G}.)'])"
$ws.Range("A3").Value = "('島', ['Basic Land — Island', '({T}: Add {U}.)'])"
$ws.Range("A4").Value = "('Jaya Ballard, Task Mage', ['{1}{R}{R}', 'Legendary Creature — Human Spellshaper', '{R}, {T}, Discard a card: Destroy target blue permanent.', '{1}{R}, {T}, Discard a card: Jaya Ballard, Task Mage deals 3 damage to any target. A creature dealt damage this way can’t be regenerated this turn.', '{5}{R}{R}, {T}, Discard a card: Jaya Ballard deals 6 damage to each creature and each player.', '2/2'])"
$ws.Range("A5").Value = "('山', ['Basic Land — Mountain', '({T}: Add {R}.)'])"
$ws.Range("A6").Value = "('平地', ['Basic Land — Plains', '({T}: Add {W}.)'])"
$ws.Range("A7").Value = "('沼', ['Basic Land — Swamp', '({T}: Add {B}.)'])"

# Remove the now-obsolete rows 8 through 23
$ws.Range("A8:A23").ClearContents()
